# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 6;   I = "sv"; J = "Statement-opinion" },
    @{ Row = 11;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 42;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 54;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 70;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 79;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 83;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 101; I = "aa"; J = "Agree/Accept" },
    @{ Row = 111; I = "ba"; J = "Appreciation" },
    @{ Row = 127; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 129; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 138; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 141; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 144; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.I
    $ws.Range("J" + $u.Row).Value = $u.J
}
